$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell "D2" '306.56'
Set-TextCell "E2" '-0.35%'
Set-TextCell "D3" '40.34'
Set-TextCell "E3" '1.01%'
Set-TextCell "E4" '0.83%'
Set-TextCell "D5" '0.07581'
Set-TextCell "E5" '-2.51%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell "D6" '1.608'
Set-TextCell "E6" '-2.60%'
$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell "D7" '2.447'
Set-TextCell "E7" '-4.41%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell "D8" '0.9063'
Set-TextCell "E8" '-1.48%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell "D9" '0.1011'
Set-TextCell "E9" '2.40%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell "D10" '0.1753'
Set-TextCell "E10" '0.93%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell "D11" '0.09086'
Set-TextCell "E11" '1.73%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell "D12" '0.04222'
Set-TextCell "E12" '-4.03%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell "D13" '0.1054'
Set-TextCell "E13" '-0.47%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell "D14" '0.001226'
Set-TextCell "E14" '-2.22%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell "D15" '0.005830'
Set-TextCell "E15" '3.23%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell "D16" '3.348'
Set-TextCell "E16" '-0.52%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell "D17" '4.274'
Set-TextCell "E17" '-1.12%'
Set-TextCell "D18" '0.3274'
Set-TextCell "E18" '-2.86%'
Set-TextCell "D19" '6.664'
Set-TextCell "E19" '-5.47%'
Set-TextCell "D20" '0.1358'
Set-TextCell "E20" '-0.40%'
Set-TextCell "E21" '2.74%'
Set-TextCell "D22" '0.04181'
Set-TextCell "E22" '0.87%'
Set-TextCell "E23" '1.95%'
Set-TextCell "D24" '0.004049'
Set-TextCell "E24" '-0.92%'
Set-TextCell "D25" '0.0001302'
Set-TextCell "E25" '6.22%'
Set-TextCell "D26" '0.0003013'
Set-TextCell "E26" '0.79%'
Set-TextCell "D38" '0.02386'
Set-TextCell "E38" '-0.14%'
Set-TextCell "D39" '0.05135'
Set-TextCell "E39" '-0.92%'
Set-TextCell "D40" '0.007769'
Set-TextCell "E40" '-2.73%'
Set-TextCell "D41" '0.1294'
Set-TextCell "E41" '-2.55%'
Set-TextCell "D42" '0.007053'
Set-TextCell "E42" '-3.04%'
Set-TextCell "D43" '0.001973'
Set-TextCell "E43" '-0.32%'
Set-TextCell "D44" '0.008455'
Set-TextCell "E44" '4.83%'
Set-TextCell "D45" '0.3316'
Set-TextCell "E45" '-0.58%'
Set-TextCell "D46" '0.00006376'
Set-TextCell "E46" '-5.04%'
Set-TextCell "E47" '-0.28%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextCell "D48" '0.004409'
Set-TextCell "E48" '7.15%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextCell "D49" '0.006785'
Set-TextCell "E49" '98.15%'
Set-TextCell "D50" '0.00002104'
Set-TextCell "E50" '-0.28%'
Set-TextCell "D51" '0.0002004'
Set-TextCell "E51" '-0.28%'
